$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Trim "Wodehouse woke in the morning having slept ..." down to
#    "Wodehouse slept ..." by deleting the phrase
#    "woke in the morning having " (trailing space included so the
#    remaining words rejoin cleanly as "Wodehouse slept").
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("woke in the morning having ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $cutStart = $findRange.Start
    $cutEnd = $findRange.End

    $cutRange = $d.Range($cutStart, $cutEnd)
    $cutRange.Delete()

    # --------------------------------------------------------------
    # 2) Word keeps a single "_GoBack" bookmark marking the location
    #    of the most recent edit. Re-adding a bookmark with that
    #    reserved name moves it (a document can only have one
    #    bookmark per name), which both plants it at the new edit
    #    site - right before "slept through the resulting
    #    commotion." - and removes it from its old location next to
    #    "Kent" earlier in the document.
    # --------------------------------------------------------------
    $goBackRange = $d.Range($cutStart, $cutStart)
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
}
